# Update the "reporting" column (Q) for all employee rows: the old reporting
# manager e-mail (chandana.vennam@optimworks.com) is replaced with the valid
# reporter mounikareddy@optimworks.com, and the cells are turned into mailto
# hyperlinks (matching the existing email-hyperlink pattern already used in
# column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newReporter = "mounikareddy@optimworks.com"
$mailto = "mailto:" + $newReporter

# 1. Update the cell text for every data row (2-7).
$ws.Range("Q2:Q7").Value = $newReporter

# 2. Re-create the hyperlinks on the Q column. Q2 gets its own hyperlink,
#    Q3:Q7 share a single hyperlink entry (mirrors how Excel collapses a
#    contiguous, identically-targeted range into one <hyperlink> element).
$ws.Hyperlinks.Add($ws.Range("Q2"), $mailto)
$ws.Hyperlinks.Add($ws.Range("Q3:Q7"), $mailto, [Type]::Missing, [Type]::Missing, $newReporter)

# 3. Apply the built-in Hyperlink style so the cells look like the other
#    mailto links already on the sheet (column D).
$ws.Range("Q2:Q7").Style = "Hyperlink"

# 4. Widen column Q so the longer address is fully visible.
$ws.Columns.Item(17).ColumnWidth = 43.72

# 5. Move the selection to reflect where the edit left the cursor.
$ws.Range("Q10").Select()
